$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$t  = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------------
# 1. Shift the data of columns H:M (rows 1-7) one column to the right, so
#    a brand-new, blank column opens up at H. Work from the rightmost
#    column to the left so a source column is never clobbered before it has
#    been read. Range.Copy carries both value AND number/visual format, so
#    the style index travels together with the cell exactly like a native
#    Excel "insert column" does.
# ---------------------------------------------------------------------------
$ws.Range("M1:M7").Copy($ws.Range("N1:N7")) | Out-Null
$ws.Range("L1:L7").Copy($ws.Range("M1:M7")) | Out-Null
$ws.Range("K1:K7").Copy($ws.Range("L1:L7")) | Out-Null
$ws.Range("J1:J7").Copy($ws.Range("K1:K7")) | Out-Null
$ws.Range("I1:I7").Copy($ws.Range("J1:J7")) | Out-Null
$ws.Range("H1:H7").Copy($ws.Range("I1:I7")) | Out-Null

# ---------------------------------------------------------------------------
# 2. Populate the freshly opened column H with the new "EventAdd" column.
#    Re-use the exact same styles the header/sub-header/blank-data cells had
#    before the shift (they are still readable from column I, which just
#    received a copy of the old H content/style).
# ---------------------------------------------------------------------------
$ws.Range("H1").Value = "增加"
$ws.Range("H1").Style = "Normal"
$ws.Range("H1").Font.Size = $ws.Range("H1").Font.Size
$ws.Range("H2").Value = "string[]"
$ws.Range("H3").Value = "EventAdd"

$ws.Range("H1:H3").Copy() | Out-Null
# Put back the styling that belongs to the header rows (copied from I, the
# column that used to be H) - values must stay put, so reapply style only.
$h1Style = $ws.Range("I1").Style
$h2Style = $ws.Range("I2").Style
$h3Style = $ws.Range("I3").Style

$ws.Range("H1").Style = $h1Style
$ws.Range("H2").Style = $h2Style
$ws.Range("H3").Style = $h3Style

$ws.Range("H4").Style = $ws.Range("I4").Style
$ws.Range("H5").Style = $ws.Range("I5").Style
$ws.Range("H6").Style = $ws.Range("I6").Style
$ws.Range("H7").Style = $ws.Range("I7").Style

# ---------------------------------------------------------------------------
# 3. Grow the table so it knows about the new column and the new row, then
#    add the new story row's data.
# ---------------------------------------------------------------------------
$t.Resize($ws.Range("A3:N8")) | Out-Null

$ws.Range("A8").Value = 47000021
$ws.Range("B8").Value = "远古之路"
$ws.Range("C8").Value = "从冰冻苔原往前走,在亚瑞特山脚下，就是传说中的远古之路。再往前就是亚瑞特山脉颠峰，在那里有3个古代勇士守护者通往远古遗迹的道路。"
$ws.Range("D8").Value = 100
$ws.Range("E8").Value = 18000201

$ws.Range("A8").Style = $ws.Range("A7").Style
$ws.Range("B8:E8").Style = $ws.Range("B7").Style

$ws.Range("F8:N8").Value = 0
$ws.Range("F8:N8").ClearContents() | Out-Null
$ws.Range("F8:N8").Font.Name = "宋体"

# ---------------------------------------------------------------------------
# 4. Column widths: a new, narrower column H ("EventAdd") is inserted and
#    the former widths slide one slot to the right.
# ---------------------------------------------------------------------------
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("M").ColumnWidth = $ws.Columns("L").ColumnWidth
$ws.Columns("L").ColumnWidth = $ws.Columns("K").ColumnWidth
$ws.Columns("K").ColumnWidth = $ws.Columns("J").ColumnWidth
$ws.Columns("J").ColumnWidth = $ws.Columns("I").ColumnWidth
$ws.Columns("I").ColumnWidth = $ws.Columns("H").ColumnWidth
$ws.Columns("H").ColumnWidth = 13.875

# ---------------------------------------------------------------------------
# 5. Conditional formatting ranges grow/move together with the data they
#    watch (same adjustment Excel performs automatically on a real column
#    insert / table resize).
# ---------------------------------------------------------------------------
$cf1 = $ws.Range("G4:J7 M4:M7").FormatConditions.Item(1)
$cf1.ModifyAppliesToRange($ws.Range("G4:K8 N4:N8"))

$cf2 = $ws.Range("K4:K7").FormatConditions.Item(1)
$cf2.ModifyAppliesToRange($ws.Range("L4:L8"))

$cf3 = $ws.Range("L4:L7").FormatConditions.Item(1)
$cf3.ModifyAppliesToRange($ws.Range("M4:M8"))

$cf4 = $ws.Range("F4:F7").FormatConditions.Item(1)
$cf4.ModifyAppliesToRange($ws.Range("F4:F8"))

# ---------------------------------------------------------------------------
# 6. Cosmetic: move the active selection like the saved workbook shows.
# ---------------------------------------------------------------------------
$ws.Range("I4").Select() | Out-Null
